$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GanttChart")
$shp = $ws.Shapes.Item("Scroll Bar 46")
$cf = $shp.ControlFormat
$cf.Value = 5
$ws.Range("H4").Value = 5
